$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4545850.5
$ws.Range("I12").Value = 15151851
$ws.Range("K12").Value = 15151851
$ws.Range("M12").Value = -15151681
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H74").Value = 4879
$ws.Range("I74").Value = 4879
$ws.Range("K74").Value = 4879
$ws.Range("M74").Value = -3943
$ws.Range("H77").Value = 4879
$ws.Range("I77").Value = 4879
$ws.Range("K77").Value = 24395
$ws.Range("M77").Value = -19715
$ws.Range("H132").Value = 36227.38
$ws.Range("I132").Value = 37414.07
$ws.Range("K132").Value = 112242.21
$ws.Range("M132").Value = -109712.21
$ws.Range("H138").Value = 4999.0493
$ws.Range("J138").Value = 3553.2322
$ws.Range("L138").Value = 10659.6966
$ws.Range("N138").Value = -20939.6966
$ws.Range("H141").Value = 1852.25
$ws.Range("I141").Value = 1566.091
$ws.Range("K141").Value = 4698.272999999999
$ws.Range("M141").Value = 481.7270000000008
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4586
$ws.Range("H42").Value = 49443.5
$ws.Range("J42").Value = 48888
$ws.Range("L42").Value = 48888
$ws.Range("N42").Value = -49860
$ws.Range("H74").Value = 545611
$ws.Range("I74").Value = 1455.6296
$ws.Range("K74").Value = 1455.6296
$ws.Range("M74").Value = -581.6296
$ws.Range("H77").Value = 545611
$ws.Range("I77").Value = 1455.6296
$ws.Range("K77").Value = 7278.148
$ws.Range("M77").Value = -2910.148
$ws.Range("H88").Value = 1734.8125
$ws.Range("I88").Value = 1715.7778
$ws.Range("J88").Value = 1759.2858
$ws.Range("K88").Value = 1715.7778
$ws.Range("L88").Value = 1759.2858
$ws.Range("M88").Value = -1309.7778
$ws.Range("N88").Value = -2571.2858
$ws.Range("H91").Value = 1734.8125
$ws.Range("I91").Value = 1715.7778
$ws.Range("J91").Value = 1759.2858
$ws.Range("K91").Value = 1715.7778
$ws.Range("L91").Value = 1759.2858
$ws.Range("M91").Value = -311.7778000000001
$ws.Range("N91").Value = -4567.2858
$ws.Range("H97").Value = 6223.5293
$ws.Range("I97").Value = 7498.643
$ws.Range("J97").Value = 273
$ws.Range("K97").Value = 7498.643
$ws.Range("L97").Value = 273
$ws.Range("M97").Value = -7002.643
$ws.Range("N97").Value = -1265
$ws.Range("H132").Value = 3266.92
$ws.Range("I132").Value = 3050.238
$ws.Range("K132").Value = 9150.714
$ws.Range("M132").Value = -6620.714
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1173.6666
$ws.Range("I20").Value = 1269.5555
$ws.Range("J20").Value = 1029.8334
$ws.Range("K20").Value = 1269.5555
$ws.Range("L20").Value = 1029.8334
$ws.Range("M20").Value = -1022.5555
$ws.Range("N20").Value = -1523.8334
$ws.Range("H134").Value = 16982848
$ws.Range("I134").Value = 1519.75
$ws.Range("J134").Value = 180003600
$ws.Range("K134").Value = 4559.25
$ws.Range("L134").Value = 540010800
$ws.Range("M134").Value = -2024.25
$ws.Range("N134").Value = -540015870
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 99998.336
$ws.Range("J87").Value = 99998.336
$ws.Range("L87").Value = 99998.336
$ws.Range("N87").Value = -102370.336
$ws.Range("H90").Value = 99998.336
$ws.Range("J90").Value = 99998.336
$ws.Range("L90").Value = 299995.008
$ws.Range("N90").Value = -311851.008
$ws.Range("H132").Value = 30042.805
$ws.Range("I132").Value = 40578.42
$ws.Range("J132").Value = 2650.2
$ws.Range("K132").Value = 121735.26
$ws.Range("L132").Value = 7950.599999999999
$ws.Range("M132").Value = -119205.26
$ws.Range("N132").Value = -13010.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 2213.861
$ws.Range("J24").Value = 2307.8235
$ws.Range("L24").Value = 6923.470499999999
$ws.Range("N24").Value = -7383.470499999999
$ws.Range("H29").Value = 420
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 1500
$ws.Range("N29").Value = -2054
$ws.Range("H41").Value = 2166.3333
$ws.Range("I41").Value = 2166.3333
$ws.Range("K41").Value = 6498.999899999999
$ws.Range("M41").Value = -6160.999899999999
$ws.Range("H42").Value = 50013350
$ws.Range("J42").Value = 50013350
$ws.Range("L42").Value = 150040050
$ws.Range("N42").Value = -150041118
$ws.Range("H46").Value = 333338300
$ws.Range("J46").Value = 500003700
$ws.Range("L46").Value = 1500011100
$ws.Range("N46").Value = -1500011282
$ws.Range("H64").Value = 7915.3335
$ws.Range("I64").Value = 1012
$ws.Range("J64").Value = 9296
$ws.Range("K64").Value = 3036
$ws.Range("L64").Value = 27888
$ws.Range("M64").Value = -2766
$ws.Range("N64").Value = -28428
$ws.Range("H67").Value = 7915.3335
$ws.Range("I67").Value = 1012
$ws.Range("J67").Value = 9296
$ws.Range("K67").Value = 3036
$ws.Range("L67").Value = 27888
$ws.Range("M67").Value = -2100
$ws.Range("N67").Value = -29760
$ws.Range("H95").Value = 18685.625
$ws.Range("I95").Value = 15125
$ws.Range("J95").Value = 22246.25
$ws.Range("K95").Value = 45375
$ws.Range("L95").Value = 66738.75
$ws.Range("M95").Value = -43316
$ws.Range("N95").Value = -70856.75
$ws.Range("H140").Value = 23812520
$ws.Range("I140").Value = 33335780
$ws.Range("K140").Value = 100007340
$ws.Range("M140").Value = -100002160
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8995
$ws.Range("J10").Value = 8995
$ws.Range("L10").Value = 8995
$ws.Range("N10").Value = -9333
$ws.Range("H15").Value = 36548.75
$ws.Range("J15").Value = 45400
$ws.Range("L15").Value = 45400
$ws.Range("N15").Value = -45976
$ws.Range("H22").Value = 4995
$ws.Range("J22").Value = 4995
$ws.Range("L22").Value = 4995
$ws.Range("N22").Value = -6053
$ws.Range("H36").Value = 4749.1665
$ws.Range("J36").Value = 5448.75
$ws.Range("L36").Value = 5448.75
$ws.Range("N36").Value = -6418.75
$ws.Range("H81").Value = 36548.75
$ws.Range("J81").Value = 45400
$ws.Range("L81").Value = 45400
$ws.Range("N81").Value = -47396
$ws.Range("H84").Value = 36548.75
$ws.Range("J84").Value = 45400
$ws.Range("L84").Value = 136200
$ws.Range("N84").Value = -146184
$ws.Range("H102").Value = 2601
$ws.Range("I102").Value = 2214.5
$ws.Range("J102").Value = 3683.2
$ws.Range("K102").Value = 2214.5
$ws.Range("L102").Value = 3683.2
$ws.Range("M102").Value = -592.5
$ws.Range("N102").Value = -6927.2
$ws.Range("H132").Value = 1351729
$ws.Range("I132").Value = 2158
$ws.Range("J132").Value = 8549441
$ws.Range("K132").Value = 6474
$ws.Range("L132").Value = 25648323
$ws.Range("M132").Value = -3944
$ws.Range("N132").Value = -25653383
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 880
$ws.Range("I31").Value = 900
$ws.Range("K31").Value = 900
$ws.Range("M31").Value = -652
$ws.Range("H61").Value = 2688.8333
$ws.Range("I61").Value = 2611.7646
$ws.Range("K61").Value = 2611.7646
$ws.Range("M61").Value = -2409.7646
$ws.Range("H68").Value = 4999.5
$ws.Range("J68").Value = 4999.5
$ws.Range("L68").Value = 4999.5
$ws.Range("N68").Value = -6497.5
$ws.Range("H71").Value = 4999.5
$ws.Range("J71").Value = 4999.5
$ws.Range("L71").Value = 24997.5
$ws.Range("N71").Value = -32485.5
$ws.Range("H93").Value = 2822.7693
$ws.Range("I93").Value = 1569.8
$ws.Range("K93").Value = 1569.8
$ws.Range("M93").Value = -321.8
$ws.Range("H113").Value = 2688.8333
$ws.Range("I113").Value = 2611.7646
$ws.Range("K113").Value = 2611.7646
$ws.Range("M113").Value = -441.7646
$ws.Range("H132").Value = 2681.8086
$ws.Range("I132").Value = 2265.6177
$ws.Range("J132").Value = 3770.3076
$ws.Range("K132").Value = 6796.853099999999
$ws.Range("L132").Value = 11310.9228
$ws.Range("M132").Value = -4266.853099999999
$ws.Range("N132").Value = -16370.9228
$ws.Range("H136").Value = 1367.35
$ws.Range("I136").Value = 1769.8096
$ws.Range("J136").Value = 1260.3671
$ws.Range("K136").Value = 5309.4288
$ws.Range("L136").Value = 3781.1013
$ws.Range("M136").Value = -2759.4288
$ws.Range("N136").Value = -8881.1013
$ws.Range("H140").Value = 79999
$ws.Range("J140").Value = 79999
$ws.Range("L140").Value = 79999
$ws.Range("N140").Value = -90359
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 8658.333000000001
$ws.Range("J23").Value = 8658.333000000001
$ws.Range("L23").Value = 8658.333000000001
$ws.Range("N23").Value = -9116.333000000001
$ws.Range("H62").Value = 2618.8572
$ws.Range("J62").Value = 1853
$ws.Range("L62").Value = 1853
$ws.Range("N62").Value = -3101
$ws.Range("H65").Value = 2618.8572
$ws.Range("J65").Value = 1853
$ws.Range("L65").Value = 9265
$ws.Range("N65").Value = -15505
$ws.Range("H132").Value = 1968.7778
$ws.Range("I132").Value = 1489.091
$ws.Range("K132").Value = 4467.272999999999
$ws.Range("M132").Value = -1937.272999999999
$ws.Range("H136").Value = 34142
$ws.Range("I136").Value = 50460.25
$ws.Range("K136").Value = 151380.75
$ws.Range("M136").Value = -148830.75
$ws.Range("N136").Value = -18517.3638
